# Fruta / hortaliza, semanal
# Insert three new weekly rows of data right after the existing row 1006,
# pushing the existing rows 1007-1064 down to 1010-1067, and fill in the
# new rows with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 1007 (existing content shifts down).
$ws.Rows("1007:1009").Insert()

# Common / constant columns for this subset (Vega Modelo de Temuco - Limón).
$marketId    = 10
$market      = "Vega Modelo de Temuco"
$region      = "La Araucanía"
$codreg      = 9
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102003
$categoria   = "Limón"

function Set-DataRow($row, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $origen, $precioKg, $kgUnidad) {
    $ws.Cells.Item($row, 1).Value2  = $marketId
    $ws.Cells.Item($row, 2).Value2  = $market
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $variedad
    $ws.Cells.Item($row, 12).Value2 = $calidad
    $ws.Cells.Item($row, 13).Value2 = $volumen
    $ws.Cells.Item($row, 14).Value2 = $precioMin
    $ws.Cells.Item($row, 15).Value2 = $precioMax
    $ws.Cells.Item($row, 16).Value2 = $precioProm
    $ws.Cells.Item($row, 17).Value2 = $unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $precioKg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}

Set-DataRow 1007 44585 "Sin especificar" "1a amarillo" 200 18000 19000 18500 "$/bandeja 15 kilos" "Provincia del Elquí" 1233 15

Set-DataRow 1008 44585 "Sin especificar" "1a plateado" 400 22000 23000 22500 "$/caja 18 kilos" "Provincia del Elquí" 1250 18

Set-DataRow 1009 44585 "Sin especificar" "2a amarillo" 5 450000 450000 450000 "$/bins (450 kilos)" "Región de O'Higgins" 1000 450
